# Generate Report for Handback
# Updates the handoff/handback timestamp cells with newly generated values.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for the first file row.
$wsOverview.Range("G2").Value = "2016-08-16 09:06:29"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# for the first file row.
$wsZhCn.Range("H2").Value = "2016-08-16 09:06:23"
$wsZhCn.Range("K2").Value = "2016-08-16 09:06:41"

# de-de sheet: "Correspond Handoff Datetime" for the first file row.
$wsDeDe.Range("H2").Value = "2016-08-16 09:06:49"
